$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45036.50694444445
$ws.Cells.Item(2, 2).Value = 5.237
$ws.Cells.Item(2, 3).Value = 1.607
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 2.478
$ws.Cells.Item(2, 6).Value = 3.053
$ws.Cells.Item(2, 7).Value = 2.519
$ws.Cells.Item(2, 8).Value = 5.488
$ws.Cells.Item(2, 9).Value = 1.554
$ws.Cells.Item(2, 10).Value = 0.9409999999999999
$ws.Cells.Item(2, 11).Value = 4.022
$ws.Cells.Item(2, 12).Value = 1.069
$ws.Cells.Item(2, 13).Value = 0.9379999999999999
$ws.Cells.Item(2, 14).Value = 0.6929999999999999
$ws.Cells.Item(2, 15).Value = 0.87
$ws.Cells.Item(2, 16).Value = 2.749
$ws.Cells.Item(2, 17).Value = 1.106
$ws.Cells.Item(2, 18).Value = 0.51
$ws.Cells.Item(2, 19).Value = 0.063
$ws.Cells.Item(2, 20).Value = 20.094
$ws.Cells.Item(2, 21).Value = 4.803
$ws.Cells.Item(2, 22).Value = 2.372
$ws.Cells.Item(2, 23).Value = 3.891
$ws.Cells.Item(2, 24).Value = 1.034
$ws.Cells.Item(2, 25).Value = 0.249
$ws.Cells.Item(2, 26).Value = 1.98
$ws.Cells.Item(2, 27).Value = 1.154
$ws.Cells.Item(2, 28).Value = 0.674
$ws.Cells.Item(2, 29).Value = 0.9399999999999999
$ws.Cells.Item(2, 30).Value = 3.129
$ws.Cells.Item(2, 31).Value = 2.866
$ws.Cells.Item(2, 32).Value = 3.128
$ws.Cells.Item(2, 33).Value = 0.446
$ws.Cells.Item(2, 34).Value = 1.43
$ws.Cells.Item(3, 1).Value = 45036.51388888889
$ws.Cells.Item(3, 2).Value = 15.104
$ws.Cells.Item(3, 3).Value = 10.562
$ws.Cells.Item(3, 4).Value = 0.304
$ws.Cells.Item(3, 5).Value = 28.97
$ws.Cells.Item(3, 6).Value = 24.63
$ws.Cells.Item(3, 7).Value = 11.444
$ws.Cells.Item(3, 8).Value = 36.27
$ws.Cells.Item(3, 9).Value = 16.173
$ws.Cells.Item(3, 10).Value = 7.491
$ws.Cells.Item(3, 11).Value = 12.204
$ws.Cells.Item(3, 12).Value = 11.758
$ws.Cells.Item(3, 13).Value = 12.166
$ws.Cells.Item(3, 14).Value = 3.514
$ws.Cells.Item(3, 15).Value = 10.397
$ws.Cells.Item(3, 16).Value = 15.935
$ws.Cells.Item(3, 17).Value = 8.768000000000001
$ws.Cells.Item(3, 18).Value = 0.489
$ws.Cells.Item(3, 19).Value = 0.318
$ws.Cells.Item(3, 20).Value = 158.189
$ws.Cells.Item(3, 21).Value = 30.287
$ws.Cells.Item(3, 22).Value = 10.398
$ws.Cells.Item(3, 23).Value = 20.717
$ws.Cells.Item(3, 24).Value = 10.466
$ws.Cells.Item(3, 25).Value = 1.43
$ws.Cells.Item(3, 26).Value = 18.6
$ws.Cells.Item(3, 27).Value = 8.798
$ws.Cells.Item(3, 28).Value = 7.456
$ws.Cells.Item(3, 29).Value = 8.877000000000001
$ws.Cells.Item(3, 30).Value = 13.388
$ws.Cells.Item(3, 31).Value = 1.133
$ws.Cells.Item(3, 32).Value = 32.037
$ws.Cells.Item(3, 33).Value = 5.443
$ws.Cells.Item(3, 34).Value = 12.233
$ws.Cells.Item(4, 1).Value = 45036.52083333334
$ws.Cells.Item(4, 2).Value = 18.994
$ws.Cells.Item(4, 3).Value = 13.825
$ws.Cells.Item(4, 4).Value = 0.469
$ws.Cells.Item(4, 5).Value = 38.827
$ws.Cells.Item(4, 6).Value = 32.579
$ws.Cells.Item(4, 7).Value = 14.723
$ws.Cells.Item(4, 8).Value = 54.741
$ws.Cells.Item(4, 9).Value = 21.616
$ws.Cells.Item(4, 10).Value = 9.903
$ws.Cells.Item(4, 11).Value = 15.384
$ws.Cells.Item(4, 12).Value = 15.711
$ws.Cells.Item(4, 13).Value = 16.364
$ws.Cells.Item(4, 14).Value = 4.595
$ws.Cells.Item(4, 15).Value = 13.957
$ws.Cells.Item(4, 16).Value = 20.75
$ws.Cells.Item(4, 17).Value = 11.646
$ws.Cells.Item(4, 18).Value = 0.443
$ws.Cells.Item(4, 19).Value = 0.446
$ws.Cells.Item(4, 20).Value = 209.884
$ws.Cells.Item(4, 21).Value = 39.933
$ws.Cells.Item(4, 22).Value = 13.41
$ws.Cells.Item(4, 23).Value = 27.147
$ws.Cells.Item(4, 24).Value = 13.968
$ws.Cells.Item(4, 25).Value = 1.873
$ws.Cells.Item(4, 26).Value = 27.023
$ws.Cells.Item(4, 27).Value = 11.608
$ws.Cells.Item(4, 28).Value = 9.99
$ws.Cells.Item(4, 29).Value = 11.825
$ws.Cells.Item(4, 30).Value = 17.173
$ws.Cells.Item(4, 31).Value = 0.716
$ws.Cells.Item(4, 32).Value = 49.435
$ws.Cells.Item(4, 33).Value = 7.317
$ws.Cells.Item(4, 34).Value = 16.258
$ws.Cells.Item(5, 1).Value = 45036.52777777778
$ws.Cells.Item(5, 2).Value = 6.3
$ws.Cells.Item(5, 3).Value = 4.43
$ws.Cells.Item(5, 4).Value = 0.07000000000000001
$ws.Cells.Item(5, 5).Value = 11.88
$ws.Cells.Item(5, 6).Value = 10.19
$ws.Cells.Item(5, 7).Value = 4.82
$ws.Cells.Item(5, 8).Value = 22.79
$ws.Cells.Item(5, 9).Value = 6.6
$ws.Cells.Item(5, 10).Value = 3.16
$ws.Cells.Item(5, 11).Value = 5.08
$ws.Cells.Item(5, 12).Value = 4.89
$ws.Cells.Item(5, 13).Value = 4.94
$ws.Cells.Item(5, 14).Value = 1.46
$ws.Cells.Item(5, 15).Value = 4.27
$ws.Cells.Item(5, 16).Value = 6.78
$ws.Cells.Item(5, 17).Value = 3.59
$ws.Cells.Item(5, 18).Value = 0.28
$ws.Cells.Item(5, 19).Value = 0.09
$ws.Cells.Item(5, 20).Value = 61.2
$ws.Cells.Item(5, 21).Value = 12.85
$ws.Cells.Item(5, 22).Value = 4.33
$ws.Cells.Item(5, 23).Value = 8.800000000000001
$ws.Cells.Item(5, 24).Value = 4.35
$ws.Cells.Item(5, 25).Value = 0.59
$ws.Cells.Item(5, 26).Value = 10.63
$ws.Cells.Item(5, 27).Value = 3.66
$ws.Cells.Item(5, 28).Value = 3.06
$ws.Cells.Item(5, 29).Value = 3.65
$ws.Cells.Item(5, 30).Value = 5.61
$ws.Cells.Item(5, 31).Value = 0.52
$ws.Cells.Item(5, 32).Value = 20.78
$ws.Cells.Item(5, 33).Value = 2.22
$ws.Cells.Item(5, 34).Value = 5.04

$ws.Rows(6).Delete()
